$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.337.36"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "2.155.11"
$ws.Range("E3").Value = "  +3.06%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.93"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.13"
$ws.Range("E7").Value = "  +4.15%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.393"
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0860"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.14"
$ws.Range("E12").Value = "  +7.69%  "
$ws.Range("D13").Value = "2.477.69"
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.33"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.817"
$ws.Range("E15").Value = "  +2.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.56"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").Value = "2.153.03"
$ws.Range("E17").Value = "  +3.01%  "
$ws.Range("D18").Value = "39.514.62"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.37"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").Value = "0.0₃0853"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.89"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.78"
$ws.Range("E26").Value = "  +3.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.38"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.138"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.62"
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.82"
$ws.Range("E34").Value = "  +2.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.16"
$ws.Range("E35").Value = "  +11.47%  "
$ws.Range("E36").Value = "  +1.29%  "
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.59"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.24"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("E41").Value = "  +2.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.07"
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("D43").Value = "1.536.21"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("E44").Value = "  +6.30%  "
$ws.Range("E45").Value = "  +6.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0925"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("D50").Value = "2.361.65"
$ws.Range("E50").Value = "  +3.26%  "
$ws.Range("E51").Value = "  -0.49%  "

# Reset number format styling on text-forced cells back to default so style index matches original
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D46").Style = "Normal"
